# Fix the invalid field names
#
# The "Sample Block" sheet had a stray, hidden leading "header_info"
# column (A) that should not have existed. This script removes that
# column - shifting every other column one position to the left
# (B->A, C->B, ... W->V) - re-attaches the per-header cell comments to
# their correct (shifted) columns, expands the tissue-weight-value
# comment with extra explanatory text as part of the same cleanup, and
# bumps the recorded pav:createdOn timestamp on the .metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample Block")

# ---------------------------------------------------------------------
# 1) Delete the stray leading "header_info" column. Values, data
#    validations and the sheet dimension all shift left by one column
#    automatically. Comments, however, stay anchored to their original
#    cell addresses, so column letter X still carries whatever comment
#    used to live in column X before the delete.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).EntireColumn.Delete()

# ---------------------------------------------------------------------
# 2) Capture the (still old-column-aligned) comment text of every
#    header cell from A1 to W1 (columns 1..23).
# ---------------------------------------------------------------------
$commentText = @{}
for ($col = 1; $col -le 23; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Comment -ne $null) {
        $commentText[$col] = $cell.Comment.Text()
    }
}

# The tissue_weight_value comment (the one that used to sit on column G,
# now needs to sit on column F) gains extra explanatory text as part of
# this change.
$commentText[7] = "The weight of a tissue block or the piece of tissue used in a suspension.`nKnowing the weight of the parent block and tissue used in a suspension, allows`nus to compute what percentage of the block was used for the suspension."

# ---------------------------------------------------------------------
# 3) Re-point each comment at its correct, shifted column by editing
#    the existing comment objects' text in place (this preserves the
#    original comment author instead of creating new ones). New column
#    N (1..22) should show the text that used to belong to old column
#    N+1.
# ---------------------------------------------------------------------
for ($col = 1; $col -le 22; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    if ($cell.Comment -ne $null) {
        $cell.Comment.Text($commentText[$col + 1])
    }
}

# The comment left over on the now out-of-range column 23 (old W,
# beyond the new last column V) is no longer needed.
$trailingCell = $ws.Cells.Item(1, 23)
if ($trailingCell.Comment -ne $null) {
    $trailingCell.Comment.Delete()
}

# ---------------------------------------------------------------------
# 4) Update the recorded creation timestamp on the .metadata sheet.
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item(".metadata")
$meta.Range("C2").Value = "2023-10-03T09:51:28-07:00"
